$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data. Cells whose new text would otherwise be
# auto-parsed as a number are briefly forced to Text format, then restored to the
# original (Normal) cell style so only the value/type changes, matching the source data.
$ws.Range("D2").Value = "37.530.94"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.067.17"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "2.372.40"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "2.055.79"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "37.441.14"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0226"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.491.32"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "2.254.89"
$ws.Range("E51").Value = "  -1.00%  "
